$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name / tab label (date rolled from 12-17 to 12-18)
$ws.Name = "Through 2022-12-18"

# Update the "December" row label (row 13, column A)
$ws.Range("A13").Value = "December (through 12-18)"

# Update December row (row 13) values
$ws.Range("B13").Value = 23
$ws.Range("C13").Value = 57
$ws.Range("D13").Value = 70
$ws.Range("E13").Value = 42
$ws.Range("F13").Value = 29
$ws.Range("G13").Value = 83
$ws.Range("H13").Value = 132
$ws.Range("I13").Value = 78

# Update Total row (row 14) values
$ws.Range("B14").Value = 314
$ws.Range("C14").Value = 620
$ws.Range("D14").Value = 891
$ws.Range("E14").Value = 724
$ws.Range("F14").Value = 563
$ws.Range("G14").Value = 1347
$ws.Range("H14").Value = 1775
$ws.Range("I14").Value = 1595
